$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.94049999999999
$ws.Range("E4").Value = 13.60419999999999
$ws.Range("E7").Value = 11.95549999999999
$ws.Range("E8").Value = 13.54039999999999
$ws.Range("B11").Value = 5.651399999999998
$ws.Range("B12").Value = 5.906199999999997
$ws.Range("E12").Value = 11.80889999999999
$ws.Range("E14").Value = 14.0597
$ws.Range("B15").Value = 5.7863
$ws.Range("E22").Value = 12.7491
